$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6..78 down to rows 7..79
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with its data
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = 'Vega Modelo de Temuco'
$ws.Range("C6").Value = 'La Araucanía'
$ws.Range("D6").Value = (Get-Date -Year 2023 -Month 6 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D6").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 'Fruta'
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = 'Otros'
$ws.Range("I6").Value = 100107001
$ws.Range("J6").Value = 'Caqui'
$ws.Range("K6").Value = 'Fuyu'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 25
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("Q6").Value = '$/bandeja 15 kilos granel'
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 15
